# Update countries & provincias Spain
# - Re-rank Nepal, Uganda and Malta ahead of their neighbours (Kirguistan/Suiza/
#   Uzbekistan, Somalia, Yemen respectively), pushing the pre-existing rows down
#   by one position, and refresh the statistics for the affected rows.
# - Refresh statistics for a handful of other countries (Alemania, Madagascar,
#   Vietnam, Taiwan) whose ranking did not change.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-CountryRow {
    param($Row, $Country, $Total, $Nuevos, $Activos, $Recuperados, $Criticos, $MuertesHoy, $Muertes)

    $ws.Cells.Item($Row, 1).Value = $Country
    $ws.Cells.Item($Row, 2).Value = $Total
    $ws.Cells.Item($Row, 3).Value = $Nuevos
    $ws.Cells.Item($Row, 4).Value = $Activos
    $ws.Cells.Item($Row, 5).Value = $Recuperados
    $ws.Cells.Item($Row, 6).Value = $Criticos
    $ws.Cells.Item($Row, 7).Value = $MuertesHoy
    $ws.Cells.Item($Row, 8).Value = $Muertes
}

# --- Alemania: refreshed totals, no re-ranking ---
Set-CountryRow 23 "Alemania" 249063 249 224600 15064 0 0 9399

# --- Nepal jumps ahead of Kirguistan / Suiza / Uzbekistan ---
Set-CountryRow 61 "Nepal"      44236 1359 25561 18404 0 14 271
Set-CountryRow 62 "Kirguistan" 44199 64   39389 3750  0 1  1060
Set-CountryRow 63 "Suiza"      43532 405  36500 5019  0 0  2013
Set-CountryRow 64 "Uzbekistan" 42903 215  40176 2391  0 5  336

# --- Madagascar: refreshed totals, no re-ranking ---
Set-CountryRow 85 "Madagascar" 15187 81 14075 914 0 1 198

# --- Uganda jumps ahead of Somalia ---
Set-CountryRow 126 "Uganda"  3353 65 1564 1754 0 2 35
Set-CountryRow 127 "Somalia" 3310 0  2481 731  0 0 98

# --- Malta jumps ahead of Yemen ---
Set-CountryRow 146 "Malta" 1984 19 1565 406 0 0 13
Set-CountryRow 147 "Yemen" 1979 0  1180 228 0 0 571

# --- Vietnam: refreshed totals, no re-ranking ---
Set-CountryRow 165 "Vietnam" 1049 3 772 242 0 0 35

# --- Taiwan: refreshed totals, no re-ranking ---
Set-CountryRow 175 "Taiwan" 490 1 471 12 0 0 7
